# Scheduled market-data refresh for the "Lich Profits" leve-crafting workbook.
# Updates currentAveragePrice(NQ/HQ) columns (H-N) with freshly pulled market
# board values for each crafting-job worksheet; all other data is untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: refreshed market prices/profits
$ws.Range("H17").Value = 813.7742
$ws.Range("J17").Value = 813.7742
$ws.Range("L17").Value = 2441.3226
$ws.Range("N17").Value = -2777.3226

# Row 111: refreshed market prices/profits
$ws.Range("I111").Value = 9969.571
$ws.Range("J111").Value = 1498.5
$ws.Range("K111").Value = 29908.713
$ws.Range("L111").Value = 4495.5
$ws.Range("M111").Value = -26841.713
$ws.Range("N111").Value = -10629.5

# Row 118: refreshed market prices/profits
$ws.Range("H118").Value = 668.8333
$ws.Range("I118").Value = 552.7
$ws.Range("K118").Value = 1658.1
$ws.Range("M118").Value = -1.100000000000136

# Row 132: refreshed market prices/profits
$ws.Range("H132").Value = 2501.4634
$ws.Range("I132").Value = 2441.525
$ws.Range("K132").Value = 7324.575000000001
$ws.Range("M132").Value = -4794.575000000001

# Row 138: refreshed market prices/profits
$ws.Range("H138").Value = 3285.0364
$ws.Range("I138").Value = 1810.4546
$ws.Range("J138").Value = 3653.682
$ws.Range("K138").Value = 5431.3638
$ws.Range("L138").Value = 10961.046
$ws.Range("M138").Value = -291.3638000000001
$ws.Range("N138").Value = -21241.046

# Row 141: refreshed market prices/profits
$ws.Range("H141").Value = 3512
$ws.Range("I141").Value = 3085
$ws.Range("K141").Value = 9255
$ws.Range("M141").Value = -4075

$ws = $wb.Worksheets.Item("ARM")
# Row 61: refreshed market prices/profits
$ws.Range("H61").Value = 1715.0541
$ws.Range("I61").Value = 1607.3429
$ws.Range("J61").Value = 3600
$ws.Range("K61").Value = 1607.3429
$ws.Range("L61").Value = 3600
$ws.Range("M61").Value = -1395.3429
$ws.Range("N61").Value = -4024

# Row 74: refreshed market prices/profits
$ws.Range("H74").Value = 40904.92
$ws.Range("I74").Value = 46313.934
$ws.Range("J74").Value = 1238.8334
$ws.Range("K74").Value = 46313.934
$ws.Range("L74").Value = 1238.8334
$ws.Range("M74").Value = -45439.934
$ws.Range("N74").Value = -2986.8334

# Row 77: refreshed market prices/profits
$ws.Range("H77").Value = 40904.92
$ws.Range("I77").Value = 46313.934
$ws.Range("J77").Value = 1238.8334
$ws.Range("K77").Value = 231569.67
$ws.Range("L77").Value = 6194.166999999999
$ws.Range("M77").Value = -227201.67
$ws.Range("N77").Value = -14930.167

# Row 110: refreshed market prices/profits
$ws.Range("H110").Value = 3168.6843
$ws.Range("I110").Value = 3814.4666
$ws.Range("J110").Value = 747
$ws.Range("K110").Value = 3814.4666
$ws.Range("L110").Value = 747
$ws.Range("M110").Value = -1769.4666
$ws.Range("N110").Value = -4837

# Row 136: refreshed market prices/profits
$ws.Range("H136").Value = 1715.0541
$ws.Range("I136").Value = 1607.3429
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 4822.028700000001
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -2272.028700000001
$ws.Range("N136").Value = -15900

$ws = $wb.Worksheets.Item("BSM")
# Row 107: refreshed market prices/profits
$ws.Range("H107").Value = 1984.25
$ws.Range("I107").Value = 1984.25
$ws.Range("K107").Value = 1984.25
$ws.Range("M107").Value = -64.25

# Row 134: refreshed market prices/profits
$ws.Range("H134").Value = 6042.84
$ws.Range("I134").Value = 6042.84
$ws.Range("K134").Value = 18128.52
$ws.Range("M134").Value = -15593.52

$ws = $wb.Worksheets.Item("CRP")
# Row 16: refreshed market prices/profits
$ws.Range("H16").Value = 4415.625
$ws.Range("I16").Value = 2969.9167
$ws.Range("J16").Value = 8752.75
$ws.Range("K16").Value = 2969.9167
$ws.Range("L16").Value = 8752.75
$ws.Range("M16").Value = -2682.9167
$ws.Range("N16").Value = -9326.75

# Row 31: refreshed market prices/profits
$ws.Range("H31").Value = 184310.2
$ws.Range("I31").Value = 229839.9
$ws.Range("K31").Value = 229839.9
$ws.Range("M31").Value = -229544.9

# Row 34: refreshed market prices/profits
$ws.Range("H34").Value = 184310.2
$ws.Range("I34").Value = 229839.9
$ws.Range("K34").Value = 229839.9
$ws.Range("M34").Value = -229637.9

# Row 107: refreshed market prices/profits
$ws.Range("H107").Value = 3831.9546
$ws.Range("I107").Value = 709.5454999999999
$ws.Range("K107").Value = 709.5454999999999
$ws.Range("M107").Value = 1210.4545

# Row 113: refreshed market prices/profits
$ws.Range("H113").Value = 4415.625
$ws.Range("I113").Value = 2969.9167
$ws.Range("J113").Value = 8752.75
$ws.Range("K113").Value = 2969.9167
$ws.Range("L113").Value = 8752.75
$ws.Range("M113").Value = -799.9167000000002
$ws.Range("N113").Value = -13092.75

# Row 132: refreshed market prices/profits
$ws.Range("H132").Value = 7346.7837
$ws.Range("I132").Value = 6613.212
$ws.Range("K132").Value = 19839.636
$ws.Range("M132").Value = -17309.636

# Row 134: refreshed market prices/profits
$ws.Range("H134").Value = 2552.4888
$ws.Range("I134").Value = 2426.6978
$ws.Range("J134").Value = 5257
$ws.Range("K134").Value = 7280.0934
$ws.Range("L134").Value = 15771
$ws.Range("M134").Value = -4745.0934
$ws.Range("N134").Value = -20841

# Row 135: refreshed market prices/profits
$ws.Range("H135").Value = 86665.664
$ws.Range("J135").Value = 99998.5
$ws.Range("L135").Value = 99998.5
$ws.Range("N135").Value = -110138.5

$ws = $wb.Worksheets.Item("CUL")
# Row 133: refreshed market prices/profits
$ws.Range("H133").Value = 13444
$ws.Range("I133").Value = 1888
$ws.Range("K133").Value = 5664
$ws.Range("M133").Value = -604

# Row 134: refreshed market prices/profits
$ws.Range("H134").Value = 7975.8
$ws.Range("I134").Value = 7975.8
$ws.Range("K134").Value = 23927.4
$ws.Range("M134").Value = -18857.4

# Row 139: refreshed market prices/profits
$ws.Range("H139").Value = 3257.0715
$ws.Range("I139").Value = 2011
$ws.Range("J139").Value = 5500
$ws.Range("K139").Value = 6033
$ws.Range("L139").Value = 16500
$ws.Range("M139").Value = -893
$ws.Range("N139").Value = -26780

$ws = $wb.Worksheets.Item("GSM")
# Row 113: refreshed market prices/profits
$ws.Range("H113").Value = 1588.8334
$ws.Range("I113").Value = 1432.8667
$ws.Range("J113").Value = 2368.6667
$ws.Range("K113").Value = 1432.8667
$ws.Range("L113").Value = 2368.6667
$ws.Range("M113").Value = 737.1333
$ws.Range("N113").Value = -6708.6667

$ws = $wb.Worksheets.Item("LTW")
# Row 93: refreshed market prices/profits
$ws.Range("H93").Value = 47620704
$ws.Range("I93").Value = 1501.4286
$ws.Range("K93").Value = 1501.4286
$ws.Range("M93").Value = -253.4286

$ws = $wb.Worksheets.Item("WVR")
# Row 107: refreshed market prices/profits
$ws.Range("H107").Value = 1869
$ws.Range("I107").Value = 3172.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 9517.5
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -7597.5
$ws.Range("N107").Value = -6840

# Row 122: refreshed market prices/profits
$ws.Range("H122").Value = 3509.9092
$ws.Range("I122").Value = 3444.4285
$ws.Range("J122").Value = 3624.5
$ws.Range("K122").Value = 10333.2855
$ws.Range("L122").Value = 10873.5
$ws.Range("M122").Value = -7883.2855
$ws.Range("N122").Value = -15773.5

# Row 132: refreshed market prices/profits
$ws.Range("H132").Value = 2378.3704
$ws.Range("I132").Value = 1968.94
$ws.Range("J132").Value = 7496.25
$ws.Range("K132").Value = 5906.82
$ws.Range("L132").Value = 22488.75
$ws.Range("M132").Value = -3376.82
$ws.Range("N132").Value = -27548.75

# Row 136: refreshed market prices/profits
$ws.Range("H136").Value = 401302.84
$ws.Range("I136").Value = 455753.7
$ws.Range("K136").Value = 1367261.1
$ws.Range("M136").Value = -1364711.1
